# Applies the "cryptos list" refresh described in the commit diff:
# updated prices / 1h volume percentages for most rows, plus a row-order
# swap between Monero (row 28) and LidoDAOToken (row 29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value. Many "Price" (column D)
# values look like plain numbers (e.g. '1.006'); Excel's COM Value setter
# auto-converts those to numeric cells, which would change the stored cell
# type away from text and reformat it (e.g. "0.00001143" -> 1.143E-05).
# Prefixing with a leading single quote forces text entry (matching the
# original workbook), and resetting .Style back to "Normal" afterwards
# clears the quote-prefix formatting so no stray style gets attached to
# the cell itself.
$updates = [ordered]@{
    'D2' = '30.460.50'
    'E2' = '  +0.05%  '
    'D3' = '2.105.75'
    'E3' = '  -0.06%  '
    'D4' = '1.006'
    'E4' = '  +0.73%  '
    'D5' = '334.72'
    'E5' = '  +1.65%  '
    'E6' = '  +0.69%  '
    'D7' = '0.5218'
    'E7' = '  -0.64%  '
    'D8' = '0.4544'
    'E8' = '  +4.45%  '
    'D9' = '54.46'
    'E9' = '  +15.65%  '
    'D10' = '0.08923'
    'E10' = '  +0.70%  '
    'D11' = '1.179'
    'E11' = '  +1.41%  '
    'D12' = '24.08'
    'E12' = '  -2.19%  '
    'D13' = '2.108.05'
    'E13' = '  +0.78%  '
    'D14' = '6.817'
    'E14' = '  +1.16%  '
    'D15' = '8.005'
    'E15' = '  +2.90%  '
    'D16' = '96.79'
    'E16' = '  +0.30%  '
    'D17' = '0.00001143'
    'E17' = '  +1.31%  '
    'D18' = '1.006'
    'E18' = '  +0.77%  '
    'D19' = '0.06644'
    'E19' = '  +0.17%  '
    'D20' = '19.19'
    'E20' = '  +1.12%  '
    'E21' = '  +0.43%  '
    'E22' = '  -0.54%  '
    'D23' = '30.528.12'
    'E23' = '  +0.11%  '
    'D24' = '12.38'
    'E24' = '  -0.06%  '
    'D25' = '2.349'
    'E25' = '  +1.01%  '
    'D26' = '2.347.57'
    'E26' = '  +0.42%  '
    'D27' = '22.15'
    'E27' = '  -1.41%  '
    'B28' = 'LidoDAOToken'
    'C28' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D28' = '2.532'
    'E28' = '  -2.46%  '
    'B29' = 'Monero'
    'C29' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D29' = '162.70'
    'E29' = '  +0.69%  '
    'D30' = '133.59'
    'E30' = '  +0.86%  '
    'D31' = '1.207'
    'E31' = '  +0.03%  '
    'E32' = '  -0.49%  '
    'E33' = '  -3.87%  '
    'D34' = '6.389'
    'E34' = '  +3.20%  '
    'D35' = '3.945'
    'E35' = '  +0.67%  '
    'D36' = '10.38'
    'E36' = '  +3.65%  '
    'D37' = '5.793'
    'E37' = '  +5.52%  '
    'D38' = '0.02578'
    'E38' = '  -0.21%  '
    'D39' = '0.06833'
    'E39' = '  +1.67%  '
    'E40' = '  +1.63%  '
    'D41' = '12.73'
    'E41' = '  +0.19%  '
    'D42' = '0.6870'
    'D43' = '1.248'
    'E43' = '  -0.64%  '
    'D44' = '2.319'
    'E44' = '  +4.47%  '
    'D45' = '14.07'
    'D46' = '0.6359'
    'E46' = '  -0.39%  '
    'E47' = '  +1.38%  '
    'D48' = '0.00000000352'
    'E48' = '  +22.71%  '
    'E49' = '  -0.27%  '
    'D50' = '83.19'
    'E50' = '  +1.25%  '
    'D51' = '1.203'
    'E51' = '  +0.28%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $cell = $ws.Range($cellRef)
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        # Numeric-looking text: force text storage via quote-prefix, then
        # drop back to the default style so the cell carries no formatting.
        $cell.Value = "'" + $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}
